$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that are permuted across rows 2-11: D, L, M, N, O, P, S
$cols = @(4, 12, 13, 14, 15, 16, 19)

# Capture the current (pre-edit) values for the rows involved, keyed by row number.
$orig = @{}
foreach ($r in 2..11) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $orig[$r] = $rowVals
}

# New row r gets the values that used to belong to row perm[r]
$perm = @{
    2  = 8
    3  = 11
    4  = 2
    5  = 3
    6  = 4
    7  = 5
    8  = 7
    9  = 6
    10 = 9
    11 = 10
}

foreach ($r in 2..11) {
    $src = $perm[$r]
    $srcVals = $orig[$src]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $srcVals[$c]
    }
}
